# Auto-generated: apply scheduled market-price/profit updates to each sheet
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3283.2632
$ws.Range("I62").Value = 3256.7
$ws.Range("J62").Value = 3312.7778
$ws.Range("K62").Value = 3256.7
$ws.Range("L62").Value = 3312.7778
$ws.Range("M62").Value = -2632.7
$ws.Range("N62").Value = -4560.7778
$ws.Range("H65").Value = 3283.2632
$ws.Range("I65").Value = 3256.7
$ws.Range("J65").Value = 3312.7778
$ws.Range("K65").Value = 16283.5
$ws.Range("L65").Value = 16563.889
$ws.Range("M65").Value = -13163.5
$ws.Range("N65").Value = -22803.889
$ws.Range("H68").Value = 35979.5
$ws.Range("J68").Value = 35979.5
$ws.Range("L68").Value = 35979.5
$ws.Range("N68").Value = -37477.5
$ws.Range("H71").Value = 35979.5
$ws.Range("J71").Value = 35979.5
$ws.Range("L71").Value = 107938.5
$ws.Range("N71").Value = -115426.5
$ws.Range("H92").Value = 2009.5883
$ws.Range("I92").Value = 483.07144
$ws.Range("J92").Value = 9133.333000000001
$ws.Range("K92").Value = 483.07144
$ws.Range("L92").Value = 9133.333000000001
$ws.Range("M92").Value = 764.9285600000001
$ws.Range("N92").Value = -11629.333
$ws.Range("H100").Value = 1460.0714
$ws.Range("I100").Value = 1487.625
$ws.Range("J100").Value = 1423.3334
$ws.Range("K100").Value = 1487.625
$ws.Range("L100").Value = 1423.3334
$ws.Range("M100").Value = -946.625
$ws.Range("N100").Value = -2505.3334
$ws.Range("H125").Value = 3635
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 3635
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 32715
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -37635
$ws.Range("H137").Value = 2523.6
$ws.Range("I137").Value = 1495.1666
$ws.Range("J137").Value = 4066.25
$ws.Range("K137").Value = 4485.4998
$ws.Range("L137").Value = 12198.75
$ws.Range("M137").Value = -1935.4998
$ws.Range("N137").Value = -17298.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1835.7778
$ws.Range("I2").Value = 1815.4
$ws.Range("J2").Value = 1937.6666
$ws.Range("K2").Value = 1815.4
$ws.Range("L2").Value = 1937.6666
$ws.Range("M2").Value = -1702.4
$ws.Range("N2").Value = -2163.6666
$ws.Range("H61").Value = 2018.6471
$ws.Range("I61").Value = 1107.4166
$ws.Range("K61").Value = 1107.4166
$ws.Range("M61").Value = -895.4166
$ws.Range("H74").Value = 1500.909
$ws.Range("I74").Value = 1618.2727
$ws.Range("K74").Value = 1618.2727
$ws.Range("M74").Value = -744.2727
$ws.Range("H77").Value = 1500.909
$ws.Range("I77").Value = 1618.2727
$ws.Range("K77").Value = 8091.363499999999
$ws.Range("M77").Value = -3723.363499999999
$ws.Range("H97").Value = 1479.7142
$ws.Range("I97").Value = 1193
$ws.Range("J97").Value = 1694.75
$ws.Range("K97").Value = 1193
$ws.Range("L97").Value = 1694.75
$ws.Range("M97").Value = -697
$ws.Range("N97").Value = -2686.75
$ws.Range("H102").Value = 1207
$ws.Range("I102").Value = 1136.25
$ws.Range("J102").Value = 1490
$ws.Range("K102").Value = 1136.25
$ws.Range("L102").Value = 1490
$ws.Range("M102").Value = 485.75
$ws.Range("N102").Value = -4734
$ws.Range("H116").Value = 1835.7778
$ws.Range("I116").Value = 1815.4
$ws.Range("J116").Value = 1937.6666
$ws.Range("K116").Value = 1815.4
$ws.Range("L116").Value = 1937.6666
$ws.Range("M116").Value = 478.5999999999999
$ws.Range("N116").Value = -6525.6666
$ws.Range("H136").Value = 2018.6471
$ws.Range("I136").Value = 1107.4166
$ws.Range("K136").Value = 3322.2498
$ws.Range("M136").Value = -772.2498000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1835.7778
$ws.Range("I3").Value = 1815.4
$ws.Range("J3").Value = 1937.6666
$ws.Range("K3").Value = 1815.4
$ws.Range("L3").Value = 1937.6666
$ws.Range("M3").Value = -1701.4
$ws.Range("N3").Value = -2165.6666
$ws.Range("H56").Value = 11970
$ws.Range("J56").Value = 11970
$ws.Range("L56").Value = 11970
$ws.Range("N56").Value = -13448
$ws.Range("H99").Value = 1714.9048
$ws.Range("I99").Value = 1109.2858
$ws.Range("K99").Value = 1109.2858
$ws.Range("M99").Value = 388.7141999999999
$ws.Range("H134").Value = 2367.5908
$ws.Range("I134").Value = 2203.65
$ws.Range("J134").Value = 4007
$ws.Range("K134").Value = 6610.950000000001
$ws.Range("L134").Value = 12021
$ws.Range("M134").Value = -4075.950000000001
$ws.Range("N134").Value = -17091

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1881.1111
$ws.Range("I31").Value = 1470.3158
$ws.Range("J31").Value = 4111.143
$ws.Range("K31").Value = 1470.3158
$ws.Range("L31").Value = 4111.143
$ws.Range("M31").Value = -1175.3158
$ws.Range("N31").Value = -4701.143
$ws.Range("H34").Value = 1881.1111
$ws.Range("I34").Value = 1470.3158
$ws.Range("J34").Value = 4111.143
$ws.Range("K34").Value = 1470.3158
$ws.Range("L34").Value = 4111.143
$ws.Range("M34").Value = -1268.3158
$ws.Range("N34").Value = -4515.143
$ws.Range("H59").Value = 26500
$ws.Range("J59").Value = 26500
$ws.Range("L59").Value = 26500
$ws.Range("N59").Value = -28790
$ws.Range("H105").Value = 1446.6666
$ws.Range("I105").Value = 1446.6666
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1446.6666
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 300.3334
$ws.Range("N105").ClearContents()
$ws.Range("H134").Value = 1860.625
$ws.Range("I134").Value = 1072.1666
$ws.Range("K134").Value = 3216.4998
$ws.Range("M134").Value = -681.4998000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 3625.8386
$ws.Range("I109").Value = 465.4
$ws.Range("J109").Value = 4233.615
$ws.Range("K109").Value = 1396.2
$ws.Range("L109").Value = 12700.845
$ws.Range("M109").Value = -356.1999999999998
$ws.Range("N109").Value = -14780.845
$ws.Range("H113").Value = 985858.3
$ws.Range("I113").Value = 2299407
$ws.Range("K113").Value = 6898221
$ws.Range("M113").Value = -6896051
$ws.Range("H131").Value = 765.9184
$ws.Range("I131").Value = 443.84616
$ws.Range("J131").Value = 882.2222
$ws.Range("K131").Value = 1331.53848
$ws.Range("L131").Value = 2646.6666
$ws.Range("M131").Value = 3708.46152
$ws.Range("N131").Value = -12726.6666
$ws.Range("H132").Value = 1011848.7
$ws.Range("J132").Value = 1011848.7
$ws.Range("L132").Value = 9106638.299999999
$ws.Range("N132").Value = -9111698.299999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 500012000
$ws.Range("J69").Value = 24000
$ws.Range("L69").Value = 24000
$ws.Range("N69").Value = -25498
$ws.Range("H70").Value = 7752.8237
$ws.Range("I70").Value = 7752.8237
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 7752.8237
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -7482.8237
$ws.Range("N70").ClearContents()
$ws.Range("H72").Value = 500012000
$ws.Range("J72").Value = 24000
$ws.Range("L72").Value = 72000
$ws.Range("N72").Value = -79488
$ws.Range("H73").Value = 7752.8237
$ws.Range("I73").Value = 7752.8237
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 7752.8237
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -6816.8237
$ws.Range("N73").ClearContents()
$ws.Range("H102").Value = 2769
$ws.Range("I102").Value = 2742.3157
$ws.Range("K102").Value = 2742.3157
$ws.Range("M102").Value = -1120.3157

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4933.3335
$ws.Range("I7").Value = 4933.3335
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 4933.3335
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -4821.3335
$ws.Range("N7").ClearContents()
$ws.Range("H40").Value = 4580.143
$ws.Range("I40").Value = 4164.143
$ws.Range("J40").Value = 5412.143
$ws.Range("K40").Value = 4164.143
$ws.Range("L40").Value = 5412.143
$ws.Range("M40").Value = -4028.143
$ws.Range("N40").Value = -5684.143
$ws.Range("H93").Value = 7076.353
$ws.Range("I93").Value = 9023.076999999999
$ws.Range("J93").Value = 749.5
$ws.Range("K93").Value = 9023.076999999999
$ws.Range("L93").Value = 749.5
$ws.Range("M93").Value = -7775.076999999999
$ws.Range("N93").Value = -3245.5
$ws.Range("H122").Value = 2290
$ws.Range("I122").Value = 1600
$ws.Range("K122").Value = 4800
$ws.Range("M122").Value = -2350
$ws.Range("H126").Value = 4933.3335
$ws.Range("I126").Value = 4933.3335
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 14800.0005
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -12330.0005
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2630.3572
$ws.Range("I122").Value = 2540
$ws.Range("K122").Value = 7620
$ws.Range("M122").Value = -5170
$ws.Range("H126").Value = 3082.5
$ws.Range("I126").Value = 3198
$ws.Range("J126").Value = 2505
$ws.Range("K126").Value = 9594
$ws.Range("L126").Value = 7515
$ws.Range("M126").Value = -7124
$ws.Range("N126").Value = -12455
$ws.Range("H132").Value = 2444.7827
$ws.Range("I132").Value = 1926.3334
$ws.Range("K132").Value = 5779.0002
$ws.Range("M132").Value = -3249.0002
$ws.Range("H136").Value = 7516.778
$ws.Range("I136").Value = 9608.5
$ws.Range("J136").Value = 3333.3333
$ws.Range("K136").Value = 28825.5
$ws.Range("L136").Value = 9999.999899999999
$ws.Range("M136").Value = -26275.5
$ws.Range("N136").Value = -15099.9999
